# New weekly price-report entry for "Pepino dulce" (Terminal La Palmera de
# La Serena) is inserted at the top of the data block (row 594), pushing the
# existing records down by one row (594-672 -> 595-673) and extending the
# sheet's used range to row 673.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row, shifting rows 594:672 down to 595:673.
$ws.Rows(594).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A594").Value = 8
$ws.Range("B594").Value = "Terminal La Palmera de La Serena"
$ws.Range("C594").Value = "Coquimbo"
$ws.Range("D594").Value = 45154
$ws.Range("E594").Value = 4
$ws.Range("F594").Value = 100112043
$ws.Range("G594").Value = "Pepino dulce"
$ws.Range("H594").Value = "Sin especificar"
$ws.Range("I594").Value = "Segunda"
$ws.Range("J594").Value = 200
$ws.Range("K594").Value = 16500
$ws.Range("L594").Value = 17000
$ws.Range("M594").Value = 16750
$ws.Range("N594").Value = "$/bandeja 18 kilos"
$ws.Range("O594").Value = "Provincia de Limarí"
$ws.Range("P594").Value = 931
$ws.Range("Q594").Value = 18
$ws.Range("R594").Value = "Hortaliza"
